$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 522.5172
$ws.Range("J17").Value = 532.25
$ws.Range("L17").Value = 1596.75
$ws.Range("N17").Value = -1932.75
$ws.Range("H19").Value = 1361.9231
$ws.Range("J19").Value = 1323.8334
$ws.Range("L19").Value = 1323.8334
$ws.Range("N19").Value = -1673.8334
$ws.Range("H41").Value = 1281.2858
$ws.Range("J41").Value = 1556.3636
$ws.Range("L41").Value = 1556.3636
$ws.Range("N41").Value = -2436.3636
$ws.Range("H106").Value = 1250.0769
$ws.Range("I106").Value = 1250.0769
$ws.Range("K106").Value = 1250.0769
$ws.Range("M106").Value = -619.0769
$ws.Range("H113").Value = 5697.5
$ws.Range("I113").Value = 3900
$ws.Range("K113").Value = 3900
$ws.Range("M113").Value = -646
$ws.Range("H120").Value = 95000
$ws.Range("J120").Value = 95000
$ws.Range("L120").Value = 95000
$ws.Range("N120").Value = -104676
$ws.Range("H129").Value = 513.1111
$ws.Range("I129").Value = 513.1111
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1539.3333
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3460.6667
$ws.Range("N129").Value = ""
$ws.Range("H135").Value = 2326.8845
$ws.Range("I135").Value = 1988.8334
$ws.Range("K135").Value = 17899.5006
$ws.Range("M135").Value = -15364.5006
$ws.Range("H138").Value = 26935.586
$ws.Range("I138").Value = 2715.6365
$ws.Range("J138").Value = 35816.234
$ws.Range("K138").Value = 8146.9095
$ws.Range("L138").Value = 107448.702
$ws.Range("M138").Value = -3006.9095
$ws.Range("N138").Value = -117728.702
$ws.Range("H141").Value = 863.4375
$ws.Range("I141").Value = 787.6667
$ws.Range("K141").Value = 2363.0001
$ws.Range("M141").Value = 2816.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18170.7
$ws.Range("I32").Value = 19064.123
$ws.Range("J32").Value = 1195.6666
$ws.Range("K32").Value = 19064.123
$ws.Range("L32").Value = 1195.6666
$ws.Range("M32").Value = -18777.123
$ws.Range("N32").Value = -1769.6666
$ws.Range("H61").Value = 8392.1875
$ws.Range("I61").Value = 1057.75
$ws.Range("K61").Value = 1057.75
$ws.Range("M61").Value = -845.75
$ws.Range("H74").Value = 873572.1
$ws.Range("I74").Value = 1500751.5
$ws.Range("K74").Value = 1500751.5
$ws.Range("M74").Value = -1499877.5
$ws.Range("H77").Value = 873572.1
$ws.Range("I77").Value = 1500751.5
$ws.Range("K77").Value = 7503757.5
$ws.Range("M77").Value = -7499389.5
$ws.Range("H102").Value = 1471.6364
$ws.Range("I102").Value = 1407.871
$ws.Range("K102").Value = 1407.871
$ws.Range("M102").Value = 214.1289999999999
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H132").Value = 1506.48
$ws.Range("J132").Value = 2788.8
$ws.Range("L132").Value = 8366.400000000001
$ws.Range("N132").Value = -13426.4
$ws.Range("H136").Value = 8392.1875
$ws.Range("I136").Value = 1057.75
$ws.Range("K136").Value = 3173.25
$ws.Range("M136").Value = -623.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1966.3334
$ws.Range("I36").Value = 1474.5
$ws.Range("J36").Value = 2950
$ws.Range("K36").Value = 1474.5
$ws.Range("L36").Value = 2950
$ws.Range("M36").Value = -940.5
$ws.Range("N36").Value = -4018
$ws.Range("H54").Value = 2491.3333
$ws.Range("I54").Value = 2491.3333
$ws.Range("K54").Value = 2491.3333
$ws.Range("M54").Value = -2007.3333
$ws.Range("H64").Value = 12647.223
$ws.Range("I64").Value = 2266.5
$ws.Range("K64").Value = 2266.5
$ws.Range("M64").Value = -2041.5
$ws.Range("H67").Value = 12647.223
$ws.Range("I67").Value = 2266.5
$ws.Range("K67").Value = 2266.5
$ws.Range("M67").Value = -1486.5
$ws.Range("H75").Value = 26453.346
$ws.Range("I75").Value = 24262.334
$ws.Range("K75").Value = 24262.334
$ws.Range("M75").Value = -23326.334
$ws.Range("H78").Value = 26453.346
$ws.Range("I78").Value = 24262.334
$ws.Range("K78").Value = 72787.00199999999
$ws.Range("M78").Value = -68107.00199999999
$ws.Range("H94").Value = 2034.9231
$ws.Range("I94").Value = 1586.5
$ws.Range("K94").Value = 1586.5
$ws.Range("M94").Value = -1135.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 69851.336
$ws.Range("J127").Value = 69851.336
$ws.Range("L127").Value = 69851.336
$ws.Range("N127").Value = -79771.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 700
$ws.Range("I8").Value = 700
$ws.Range("K8").Value = 2100
$ws.Range("M8").Value = -1961
$ws.Range("H76").Value = 4928.5713
$ws.Range("H79").Value = 4928.5713
$ws.Range("H107").Value = 3446.6365
$ws.Range("I107").Value = 30000
$ws.Range("J107").Value = 791.3
$ws.Range("K107").Value = 90000
$ws.Range("L107").Value = 2373.9
$ws.Range("M107").Value = -88080
$ws.Range("N107").Value = -6213.9
$ws.Range("H124").Value = 1950
$ws.Range("I124").Value = 1950
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 5850
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -940
$ws.Range("N124").Value = ""
$ws.Range("H140").Value = 3034.7334
$ws.Range("I140").Value = 3034.7334
$ws.Range("K140").Value = 9104.200199999999
$ws.Range("M140").Value = -3924.200199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23833.334
$ws.Range("J15").Value = 23833.334
$ws.Range("L15").Value = 23833.334
$ws.Range("N15").Value = -24409.334
$ws.Range("H42").Value = 35000
$ws.Range("J42").Value = 35000
$ws.Range("L42").Value = 35000
$ws.Range("N42").Value = -35970
$ws.Range("H57").Value = 43750
$ws.Range("J57").Value = 43750
$ws.Range("L57").Value = 43750
$ws.Range("N57").Value = -45390
$ws.Range("H70").Value = 5600.1055
$ws.Range("J70").Value = 5967.6
$ws.Range("L70").Value = 5967.6
$ws.Range("N70").Value = -6507.6
$ws.Range("H73").Value = 5600.1055
$ws.Range("J73").Value = 5967.6
$ws.Range("L73").Value = 5967.6
$ws.Range("N73").Value = -7839.6
$ws.Range("H80").Value = 8890.522999999999
$ws.Range("I80").Value = 4182.5835
$ws.Range("J80").Value = 15167.777
$ws.Range("K80").Value = 4182.5835
$ws.Range("L80").Value = 15167.777
$ws.Range("M80").Value = -3184.5835
$ws.Range("N80").Value = -17163.777
$ws.Range("H81").Value = 23833.334
$ws.Range("J81").Value = 23833.334
$ws.Range("L81").Value = 23833.334
$ws.Range("N81").Value = -25829.334
$ws.Range("H83").Value = 8890.522999999999
$ws.Range("I83").Value = 4182.5835
$ws.Range("J83").Value = 15167.777
$ws.Range("K83").Value = 20912.9175
$ws.Range("L83").Value = 75838.88499999999
$ws.Range("M83").Value = -15920.9175
$ws.Range("N83").Value = -85822.88499999999
$ws.Range("H84").Value = 23833.334
$ws.Range("J84").Value = 23833.334
$ws.Range("L84").Value = 71500.00199999999
$ws.Range("N84").Value = -81484.00199999999
$ws.Range("H97").Value = 1043.4
$ws.Range("I97").Value = 1043.4
$ws.Range("K97").Value = 1043.4
$ws.Range("M97").Value = -547.4000000000001
$ws.Range("H105").Value = 10000671
$ws.Range("J105").Value = 10000671
$ws.Range("L105").Value = 10000671
$ws.Range("N105").Value = -10007659
$ws.Range("H113").Value = 2589.2856
$ws.Range("I113").Value = 1567
$ws.Range("J113").Value = 6337.6665
$ws.Range("K113").Value = 1567
$ws.Range("L113").Value = 6337.6665
$ws.Range("M113").Value = 603
$ws.Range("N113").Value = -10677.6665
$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -37350
$ws.Range("H122").Value = 3746.7334
$ws.Range("I122").Value = 3361.6155
$ws.Range("K122").Value = 10084.8465
$ws.Range("M122").Value = -7634.8465
$ws.Range("H123").Value = 41562.312
$ws.Range("J123").Value = 41562.312
$ws.Range("L123").Value = 41562.312
$ws.Range("N123").Value = -46462.312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2576.2222
$ws.Range("J46").Value = 3106.2632
$ws.Range("L46").Value = 3106.2632
$ws.Range("N46").Value = -3482.2632
$ws.Range("H80").Value = 45000
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -28877
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 45000
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -84384
$ws.Range("N83").Value = -161232
$ws.Range("H136").Value = 3269.5293
$ws.Range("I136").Value = 2932.7896
$ws.Range("K136").Value = 8798.3688
$ws.Range("M136").Value = -6248.3688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 21374.75
$ws.Range("I81").Value = 26833
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 53666
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -52605
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 21374.75
$ws.Range("I84").Value = 26833
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 268330
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -263026
$ws.Range("N84").Value = -60608

Write-Output "applied"